$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at 413-415; existing rows 413:435 shift down to 416:438
$ws.Rows("413:415").Insert()

# Row 413
$ws.Range("A413").Value = 10
$ws.Range("B413").Value = "Vega Modelo de Temuco"
$ws.Range("C413").Value = "La Araucanía"
$ws.Range("D413").Value = 44516
$ws.Range("E413").Value = 9
$ws.Range("F413").Value = "Fruta"
$ws.Range("G413").Value = 100102
$ws.Range("H413").Value = "Cítricos"
$ws.Range("I413").Value = 100102004
$ws.Range("J413").Value = "Mandarina"
$ws.Range("K413").Value = "Murcott"
$ws.Range("L413").Value = "Especial"
$ws.Range("M413").Value = 65
$ws.Range("N413").Value = 15000
$ws.Range("O413").Value = 15000
$ws.Range("P413").Value = 15000
$ws.Range("Q413").Value = "`$/bandeja 18 kilos"
$ws.Range("R413").Value = "Región de O'Higgins"
$ws.Range("S413").Value = 833
$ws.Range("T413").Value = 18

# Row 414
$ws.Range("A414").Value = 10
$ws.Range("B414").Value = "Vega Modelo de Temuco"
$ws.Range("C414").Value = "La Araucanía"
$ws.Range("D414").Value = 44516
$ws.Range("E414").Value = 9
$ws.Range("F414").Value = "Fruta"
$ws.Range("G414").Value = 100102
$ws.Range("H414").Value = "Cítricos"
$ws.Range("I414").Value = 100102004
$ws.Range("J414").Value = "Mandarina"
$ws.Range("K414").Value = "Murcott"
$ws.Range("L414").Value = "Primera"
$ws.Range("M414").Value = 165
$ws.Range("N414").Value = 6000
$ws.Range("O414").Value = 7000
$ws.Range("P414").Value = 6394
$ws.Range("Q414").Value = "`$/bandeja 10 kilos"
$ws.Range("R414").Value = "Provincia del Elquí"
$ws.Range("S414").Value = 639
$ws.Range("T414").Value = 10

# Row 415
$ws.Range("A415").Value = 10
$ws.Range("B415").Value = "Vega Modelo de Temuco"
$ws.Range("C415").Value = "La Araucanía"
$ws.Range("D415").Value = 44516
$ws.Range("E415").Value = 9
$ws.Range("F415").Value = "Fruta"
$ws.Range("G415").Value = 100102
$ws.Range("H415").Value = "Cítricos"
$ws.Range("I415").Value = 100102004
$ws.Range("J415").Value = "Mandarina"
$ws.Range("K415").Value = "Murcott"
$ws.Range("L415").Value = "Primera"
$ws.Range("M415").Value = 155
$ws.Range("N415").Value = 10000
$ws.Range("O415").Value = 10000
$ws.Range("P415").Value = 10000
$ws.Range("Q415").Value = "`$/bandeja 18 kilos"
$ws.Range("R415").Value = "Región de O'Higgins"
$ws.Range("S415").Value = 556
$ws.Range("T415").Value = 18
